$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1 - copy the formatting from an existing header cell (H1)
# so they share the same bold/border/centered style used by the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new columns I (I0) and J (IF), rows 2-6
$ws.Cells.Item(2, 9).Value = 9
$ws.Cells.Item(2, 10).Value = 9

$ws.Cells.Item(3, 9).Value = 7
$ws.Cells.Item(3, 10).Value = 8

$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 10).Value = 9

$ws.Cells.Item(5, 9).Value = 8
$ws.Cells.Item(5, 10).Value = 9

$ws.Cells.Item(6, 9).Value = 7
$ws.Cells.Item(6, 10).Value = 8
